# Apply the "Hjemme passive updated meanEMG legmaxROM" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON)
$ws.Range("B2").Value = 2.1634793428910424
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.8658885301787187
$ws.Range("E2").ClearContents()

# Row 3 (STR)
$ws.Range("B3").Value = 2.0798837560291341
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 2.4270978639903511
$ws.Range("E3").Value = -1.2502700928198058

# Update the selection to match the new used range for the changed columns
$ws.Range("B1:E3").Select() | Out-Null
